# GPLIM-4712 add changes from discussion with users including new columns
# in Pooled Tube spreadsheet.
#
# Adds a new "Data Analysis Type" column (S) to the Pooled Tube Registration
# sheet, with a red header cell and "HybridSelection.Resequencing" filled in
# for the two sample data rows, plus two blank formatted rows below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column width to match the other bestFit columns (24.1640625 chars in
# the saved OOXML maps back to this ColumnWidth through the host's own
# char<->pixel rounding).
$ws.Columns.Item(19).ColumnWidth = 23.330729166666668

# Header cell.
$header = $ws.Range("S1")
$header.Value = "Data Analysis Type"
$header.Style = "Normal"
$header.Font.Name = "Arial"
$header.Font.Size = 10
$header.Font.ColorIndex = -4105
$header.Interior.Color = 421581

# Data rows.
$dataStyle = @($ws.Range("S2"), $ws.Range("S3"), $ws.Range("S4"), $ws.Range("S5"))
foreach ($cell in $dataStyle) {
    $cell.Style = "Normal"
    $cell.Font.Name = "Arial"
    $cell.Font.Size = 10
    $cell.Font.ColorIndex = -4105
}

$ws.Range("S2").Value = "HybridSelection.Resequencing"
$ws.Range("S3").Value = "HybridSelection.Resequencing"

Write-Host "Added Data Analysis Type column"
